# 20170601 Route 1 Fixes.
# Each timepoint (stop) was duplicated across 5 rows; trim each block back
# down to 4 rows by removing the extra (last) row of every 5-row group.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(38, 33, 28, 23, 18, 13, 8)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
